$d = $word.ActiveDocument

$pairs = @(
    @("342÷2=171, 0", "374÷6=62, 2"),
    @("859÷3=286, 1", "183÷6=30, 3"),
    @("762÷3=254, 0", "481÷4=120, 1"),
    @("682÷8=85, 2", "469÷5=93, 4"),
    @("705÷5=141, 0", "215÷3=71, 2"),
    @("631÷7=90, 1", "761÷5=152, 1"),
    @("585÷8=73, 1", "842÷5=168, 2"),
    @("428÷9=47, 5", "316÷6=52, 4"),
    @("728÷4=182, 0", "258÷6=43, 0"),
    @("869÷2=434, 1", "954÷2=477, 0"),
    @("182÷2=91, 0", "418÷8=52, 2"),
    @("156÷7=22, 2", "503÷7=71, 6"),
    @("340÷9=37, 7", "191÷4=47, 3"),
    @("185÷6=30, 5", "392÷8=49, 0"),
    @("539÷5=107, 4", "827÷2=413, 1"),
    @("574÷9=63, 7", "206÷6=34, 2"),
    @("145÷6=24, 1", "306÷8=38, 2"),
    @("622÷2=311, 0", "431÷7=61, 4"),
    @("472÷4=118, 0", "331÷6=55, 1"),
    @("458÷2=229, 0", "135÷4=33, 3"),
    @("771÷6=128, 3", "131÷3=43, 2"),
    @("711÷6=118, 3", "924÷8=115, 4"),
    @("941÷4=235, 1", "318÷3=106, 0"),
    @("973÷2=486, 1", "774÷9=86, 0"),
    @("607÷3=202, 1", "451÷7=64, 3")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
